$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (A) holds plain text values like "28/07/2022".
# Re-typing them with dashes ("28-07-2022") would make Excel silently
# reinterpret ambiguous day<=12/month<=12 strings as real dates, so force
# the cell to Text first, write the literal string, then drop the
# temporary number format again (matches the source cells, which carry no
# explicit style).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("A3")  "28-07-2022"
Set-TextValue $ws.Range("A4")  "01-08-2022"
Set-TextValue $ws.Range("A5")  "04-08-2022"
Set-TextValue $ws.Range("A6")  "08-08-2022"
Set-TextValue $ws.Range("A7")  "11-08-2022"
Set-TextValue $ws.Range("A8")  "15-08-2022"
Set-TextValue $ws.Range("A9")  "18-08-2022"
Set-TextValue $ws.Range("A10") "22-08-2022"
Set-TextValue $ws.Range("A11") "25-08-2022"
Set-TextValue $ws.Range("A12") "29-08-2022"
Set-TextValue $ws.Range("A13") "01-09-2022"
Set-TextValue $ws.Range("A14") "05-09-2022"
Set-TextValue $ws.Range("A15") "08-09-2022"
Set-TextValue $ws.Range("A16") "12-09-2022"
Set-TextValue $ws.Range("A17") "15-09-2022"
Set-TextValue $ws.Range("A18") "19-09-2022"
Set-TextValue $ws.Range("A19") "22-09-2022"
Set-TextValue $ws.Range("A20") "26-09-2022"
Set-TextValue $ws.Range("A21") "29-09-2022"

# Update attendance counters to match the revised Real/Duplicate/Invalid/Absent tally
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

$ws.Range("D11").Value = 1
$ws.Range("G11").Value = 1

$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0
